$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D values (rows 1-110) per updated amounts
$ws.Range("D1").Value = "'75070.11"
$ws.Range("D2").Value = "'24.83"
$ws.Range("D3").Value = "'120.00"
$ws.Range("D4").Value = "'755.34"
$ws.Range("D5").Value = "'41.75"
$ws.Range("D6").Value = "'1326.01"
$ws.Range("D7").Value = "'145200.00"
$ws.Range("D8").Value = "'338800.00"
$ws.Range("D9").Value = "'13112.92"
$ws.Range("D10").Value = "'0.78"
$ws.Range("D11").Value = "'8991.53"
$ws.Range("D12").Value = "'3206.86"
$ws.Range("D13").Value = "'3206.86"
$ws.Range("D14").Value = "'3206.86"
$ws.Range("D15").Value = "'10177.31"
$ws.Range("D16").Value = "'7058.11"
$ws.Range("D17").Value = "'4200.00"
$ws.Range("D18").Value = "'7485.83"
$ws.Range("D19").Value = "'400.00"
$ws.Range("D20").Value = "'3000.00"
$ws.Range("D21").Value = "'2500.00"
$ws.Range("D22").Value = "'10591.88"
$ws.Range("D23").Value = "'46432.80"
$ws.Range("D24").Value = "'1442.56"
$ws.Range("D25").Value = "'760.86"
$ws.Range("D26").Value = "'554.93"
$ws.Range("D27").Value = "'714.72"
$ws.Range("D28").Value = "'656.00"
$ws.Range("D29").Value = "'128.00"
$ws.Range("D30").Value = "'118.00"
$ws.Range("D31").Value = "'531772.00"
$ws.Range("D32").Value = "'14907.56"
$ws.Range("D33").Value = "'21.20"
$ws.Range("D34").Value = "'3879.61"
$ws.Range("D35").Value = "'99731.03"
$ws.Range("D36").Value = "'280.00"
$ws.Range("D37").Value = "'280.00"
$ws.Range("D38").Value = "'6983.23"
$ws.Range("D39").Value = "'16535.00"
$ws.Range("D40").Value = "'2775.00"
$ws.Range("D41").Value = "'20.28"
$ws.Range("D42").Value = "'439.20"
$ws.Range("D43").Value = "'478.50"
$ws.Range("D44").Value = "'35.06"
$ws.Range("D45").Value = "'376.73"
$ws.Range("D46").Value = "'440.71"
$ws.Range("D47").Value = "'264.36"
$ws.Range("D48").Value = "'7222.30"
$ws.Range("D49").Value = "'10000.00"
$ws.Range("D50").Value = "'1800.00"
$ws.Range("D51").Value = "'7512.00"
$ws.Range("D52").Value = "'2304.00"
$ws.Range("D53").Value = "'720.00"
$ws.Range("D54").Value = "'1752.00"
$ws.Range("D55").Value = "'13128.00"
$ws.Range("D56").Value = "'864.00"
$ws.Range("D57").Value = "'1440.00"
$ws.Range("D58").Value = "'3192.00"
$ws.Range("D59").Value = "'34250.00"
$ws.Range("D60").Value = "'17750.00"
$ws.Range("D61").Value = "'2500.00"
$ws.Range("D62").Value = "'2500.00"
$ws.Range("D63").Value = "'6000.00"
$ws.Range("D64").Value = "'1857.39"
$ws.Range("D65").Value = "'245.00"
$ws.Range("D66").Value = "'55.00"
$ws.Range("D67").Value = "'200.00"
$ws.Range("D68").Value = "'1500.00"
$ws.Range("D69").Value = "'15000.00"
$ws.Range("D70").Value = "'2360.00"
$ws.Range("D71").Value = "'250.00"
$ws.Range("D72").Value = "'7417.02"
$ws.Range("D73").Value = "'1800.00"
$ws.Range("D74").Value = "'40.00"
$ws.Range("D75").Value = "'2240.00"
$ws.Range("D76").Value = "'582041.25"
$ws.Range("D77").Value = "'296340.00"
$ws.Range("D78").Value = "'197560.00"
$ws.Range("D79").Value = "'246950.00"
$ws.Range("D80").Value = "'123475.00"
$ws.Range("D81").Value = "'18101.39"
$ws.Range("D82").Value = "'123475.00"
$ws.Range("D83").Value = "'48999.87"
$ws.Range("D84").Value = "'49059.96"
$ws.Range("D85").Value = "'48755.89"
$ws.Range("D86").Value = "'48865.42"
$ws.Range("D87").Value = "'1185360.00"
$ws.Range("D88").Value = "'48627.30"
$ws.Range("D89").Value = "'49191.40"
$ws.Range("D90").Value = "'50202.75"
$ws.Range("D91").Value = "'49720.00"
$ws.Range("D92").Value = "'48865.53"
$ws.Range("D93").Value = "'690.00"
$ws.Range("D94").Value = "'360.00"
$ws.Range("D95").Value = "'375.64"
$ws.Range("D96").Value = "'1000.00"
$ws.Range("D97").Value = "'2003.00"
$ws.Range("D98").Value = "'5253.00"
$ws.Range("D99").Value = "'1198.00"
$ws.Range("D100").Value = "'3000.00"
$ws.Range("D101").Value = "'2000.00"
$ws.Range("D102").Value = "'1400.00"
$ws.Range("D103").Value = "'1400.00"
$ws.Range("D104").Value = "'1300.00"
$ws.Range("D105").Value = "'150.00"
$ws.Range("D106").Value = "'700.00"
$ws.Range("D107").Value = "'5535.00"
$ws.Range("D108").Value = "'1300.00"
$ws.Range("D109").Value = "'189.13"
$ws.Range("D110").Value = "'15000.00"

# Remove the now-obsolete placeholder rows 111-114
$ws.Range("A111:FA114").EntireRow.Delete()
